$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape 1: "Content Placeholder 2" - rewrite bullet list with new structure / levels
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$tr.Text = "In short:`r6 research partners;`rhighly heterogeneous data;`rdifferent goals;`rMain challenges:`rGovernance issues;`rIntegration issues."

$levels = @(1,2,2,2,1,2,2)
$count = $tr.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.IndentLevel = $levels[$i-1]
}

# Shape 3: "TextBox 3" - update caption text
$capShape = $s.Shapes.Item(3)
$capShape.TextFrame.TextRange.Text = "Agritech project"
